$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) mixes plain decimals with source-feed values that use a second
# "." as a thousands separator (e.g. 29.688.77). Pre-format the touched cells as
# Text so Excel stores the literal digit string instead of silently parsing it as
# a number (which would also truncate trailing zeros such as "53.20" -> 53.2).
$dRows = @(2, 3, 4, 5, 7, 8, 9, 10, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26, 27, 28, 29, 30, 31, 33, 34, 35, 36, 38, 39, 40, 42, 43, 45, 46, 47, 48, 49, 50, 51)
foreach ($r in $dRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

$ws.Range("D2").Value = "29.688.77"
$ws.Range("E2").Value = "  -3.19%  "
$ws.Range("D3").Value = "2.096.18"
$ws.Range("E3").Value = "  -1.29%  "
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  -0.87%  "
$ws.Range("D5").Value = "343.25"
$ws.Range("E5").Value = "  -2.72%  "
$ws.Range("D7").Value = "0.5136"
$ws.Range("E7").Value = "  -2.82%  "
$ws.Range("D8").Value = "0.4402"
$ws.Range("E8").Value = "  -2.99%  "
$ws.Range("D9").Value = "53.20"
$ws.Range("E9").Value = "  -1.62%  "
$ws.Range("D10").Value = "0.09174"
$ws.Range("E10").Value = "  +0.88%  "
$ws.Range("E11").Value = "  -1.16%  "
$ws.Range("D12").Value = "24.90"
$ws.Range("E12").Value = "  +1.09%  "
$ws.Range("D13").Value = "2.097.57"
$ws.Range("E13").Value = "  -1.47%  "
$ws.Range("D14").Value = "6.749"
$ws.Range("E14").Value = "  -1.53%  "
$ws.Range("D15").Value = "8.148"
$ws.Range("E15").Value = "  +0.73%  "
$ws.Range("D16").Value = "99.23"
$ws.Range("E16").Value = "  -3.31%  "
$ws.Range("D17").Value = "0.00001147"
$ws.Range("E17").Value = "  -2.57%  "
$ws.Range("D18").Value = "1.009"
$ws.Range("E18").Value = "  -0.54%  "
$ws.Range("D19").Value = "20.82"
$ws.Range("E19").Value = "  +6.97%  "
$ws.Range("D20").Value = "0.06643"
$ws.Range("E20").Value = "  -1.22%  "
$ws.Range("D21").Value = "1.006"
$ws.Range("E21").Value = "  -0.52%  "
$ws.Range("D22").Value = "6.176"
$ws.Range("E22").Value = "  -2.58%  "
$ws.Range("D23").Value = "29.734.99"
$ws.Range("E23").Value = "  -3.23%  "
$ws.Range("D24").Value = "12.55"
$ws.Range("E24").Value = "  -2.15%  "
$ws.Range("D25").Value = "2.292"
$ws.Range("E25").Value = "  -4.31%  "
$ws.Range("D26").Value = "2.336.30"
$ws.Range("E26").Value = "  -1.69%  "
$ws.Range("D27").Value = "21.80"
$ws.Range("E27").Value = "  -3.13%  "
$ws.Range("D28").Value = "162.66"
$ws.Range("E28").Value = "  -1.44%  "
$ws.Range("D29").Value = "2.510"
$ws.Range("E29").Value = "  -2.21%  "
$ws.Range("D30").Value = "132.47"
$ws.Range("E30").Value = "  -2.84%  "
$ws.Range("D31").Value = "1.130"
$ws.Range("E31").Value = "  -5.83%  "
$ws.Range("E32").Value = "  -3.18%  "
$ws.Range("D33").Value = "1.632"
$ws.Range("E33").Value = "  -1.46%  "
$ws.Range("D34").Value = "6.148"
$ws.Range("E34").Value = "  -3.49%  "
$ws.Range("D35").Value = "3.968"
$ws.Range("E35").Value = "  -1.39%  "
$ws.Range("D36").Value = "6.064"
$ws.Range("E36").Value = "  -0.25%  "
$ws.Range("E37").Value = "  -1.97%  "
$ws.Range("D38").Value = "0.02563"
$ws.Range("E38").Value = "  -3.69%  "
$ws.Range("D39").Value = "0.06698"
$ws.Range("D40").Value = "12.37"
$ws.Range("E40").Value = "  -1.51%  "
$ws.Range("E41").Value = "  -1.53%  "
$ws.Range("D42").Value = "0.2214"
$ws.Range("E42").Value = "  -4.62%  "
$ws.Range("D43").Value = "1.294"
$ws.Range("E43").Value = "  +1.01%  "
$ws.Range("E44").Value = "  +2.55%  "
$ws.Range("D45").Value = "14.15"
$ws.Range("E45").Value = "  -4.85%  "
$ws.Range("D46").Value = "2.293"
$ws.Range("E46").Value = "  -1.86%  "
$ws.Range("D47").Value = "3.605"
$ws.Range("E47").Value = "  -4.36%  "
$ws.Range("D48").Value = "0.00000000344"
$ws.Range("E48").Value = "  -5.90%  "
$ws.Range("D49").Value = "1.218"
$ws.Range("E49").Value = "  -3.16%  "
$ws.Range("D50").Value = "81.79"
$ws.Range("E50").Value = "  -1.37%  "
$ws.Range("D51").Value = "1.160"
$ws.Range("E51").Value = "  -2.94%  "
